# Updates cryptos list values (Price / Volume(1h) / swapped rows 44-45)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.312.32"
$ws.Range("E2").Value = "  +4.07%  "

$ws.Range("D3").Value = "1.731.42"
$ws.Range("E3").Value = "  +2.64%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'219.38"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "'24.11"
$ws.Range("E8").Value = "  +4.78%  "

$ws.Range("E9").Value = "  +2.62%  "

$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").Value = "1.977.05"
$ws.Range("E12").Value = "  +2.70%  "

$ws.Range("D13").Value = "1.734.01"
$ws.Range("E13").Value = "  +2.67%  "

$ws.Range("D14").Value = "'4.27"
$ws.Range("E14").Value = "  +1.94%  "

$ws.Range("E15").Value = "  +1.99%  "

$ws.Range("D16").Value = "'67.85"
$ws.Range("E16").Value = "  +0.90%  "

$ws.Range("D17").Value = "28.314.81"
$ws.Range("E17").Value = "  +4.08%  "

$ws.Range("D18").Value = "'248.09"
$ws.Range("E18").Value = "  +4.12%  "

$ws.Range("E19").Value = "  +1.30%  "

$ws.Range("D20").Value = "'7.93"
$ws.Range("E20").Value = "  -2.92%  "

$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("E22").Value = "  +1.94%  "

$ws.Range("D23").Value = "'9.70"
$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("E24").Value = "  -0.33%  "

$ws.Range("D25").Value = "'149.43"
$ws.Range("E25").Value = "  +0.82%  "

$ws.Range("D26").Value = "'7.53"
$ws.Range("E26").Value = "  +3.08%  "

$ws.Range("D27").Value = "'16.72"
$ws.Range("E27").Value = "  +1.36%  "

$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").Value = "'0.0516"
$ws.Range("E30").Value = "  +2.80%  "

$ws.Range("E31").Value = "  +2.76%  "

$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("D33").Value = "'3.28"
$ws.Range("E33").Value = "  +1.17%  "

$ws.Range("D34").Value = "1.487.04"
$ws.Range("E34").Value = "  -5.55%  "

$ws.Range("E35").Value = "  -1.81%  "

$ws.Range("D36").Value = "'0.981"
$ws.Range("E36").Value = "  +2.43%  "

$ws.Range("D37").Value = "'0.602"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").Value = "'2.40"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("E39").Value = "  +1.02%  "

$ws.Range("E40").Value = "  +0.61%  "

$ws.Range("D41").Value = "'70.24"
$ws.Range("E41").Value = "  +0.93%  "

$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").Value = "'5.66"
$ws.Range("E43").Value = "  -0.71%  "

$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.880.09"
$ws.Range("E44").Value = "  +2.42%  "

$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "'2.29"
$ws.Range("E45").Value = "  +1.51%  "

$ws.Range("E46").Value = "  +1.32%  "

$ws.Range("E47").Value = "  +7.57%  "

$ws.Range("E48").Value = "  +5.15%  "

$ws.Range("D49").Value = "'90.68"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("D50").Value = "'8.22"
$ws.Range("E50").Value = "  +0.21%  "

$ws.Range("E51").Value = "  -0.81%  "
